$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 19-20 (George Floyd murdered, Protests begin)
$ws.Rows("19:20").Insert()
# Insert one new row at 22 (National Guard Activated) - after shifted row 21 (20-38)
$ws.Rows("22:22").Insert()

# Fill in new row 19: George Floyd murdered
$ws.Range("A19").Value = 43976
$ws.Range("C19").Value = "George Floyd murdered"

# Fill in new row 20: Protests begin
$ws.Range("A20").Value = 43977
$ws.Range("C20").Value = "Protests begin"

# Fill in new row 22: National Guard Activated
$ws.Range("A22").Value = 43979
$ws.Range("C22").Value = "National Guard Activated"

# Update sheet view: top-left cell and selection
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C22").Select()

# Update workbook window position
$excel.ActiveWindow.Left = 760
$excel.ActiveWindow.Top = 760
